$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Title line: "Ich hab schonmal" -> "Wer hat schonmal...?"
Replace-Text " – Ich hab schonmal " " – Wer hat schonmal…?"

# Zielbestimmung section
Replace-Text "Muss-Kriterien: Es muss local funktionieren." "Muss-Kriterien: Es muss lokal funktionieren"
Replace-Text "Wunschkriterien: siehe Lastenheft." "Wunschkriterien: siehe Lastenheft"
Replace-Text "Abgrenzungskriterien: Werbung und es darf nicht zu komplex sein." "Abgrenzungskriterien: Werbung und es darf nicht zu komplex sein"

# Produkt-Einsatz section
Replace-Text "Trinkspiel ab 16 Jahren." "Trinkspiel ab 16 Jahren"
Replace-Text "Muss zu jeder Zeit und am jedem Ort funktionieren." "Muss zu jeder Zeit und am jedem Ort funktionieren"
Replace-Text "Kann von jeder Person bedient werden." "Kann von jeder Person bedient werden"

# Produkt-Umgebung section
Replace-Text "Für offline braucht man nur ein Gerät." "Für offline braucht man nur ein Gerät"
Replace-Text "Für online einen Host, mehrere Clients und gemeinsames Netzwerk." "Für online einen Host, mehrere Clients und gemeinsames Netzwerk"

# Produkt-Funktionen section
Replace-Text "Siehe Lastenheft." "Siehe Lastenheft"

# Produkt-Daten section
Replace-Text "es die Spieldaten von nur einem Gerät." "es die Spieldaten von nur einem Gerät"
Replace-Text "und Gerätename." "und Gerätename"

# Produkt-Leistungen section
Replace-Text "Keine Begrenzung bei Nutzerzahl." "Keine Begrenzung bei Nutzerzahl"
Replace-Text "Ressourcenverbrauch der App ist gering." "Ressourcenverbrauch der App ist gering"
Replace-Text "Antwortzeit ist abhängig vom Spiel (z.B. Clients warten auf Host zum Spielstart). " "Antwortzeit ist abhängig vom Spiel (z.B. Clients warten auf Host zum Spielstart)"

# Qualitäts-Zielbestimmung section
Replace-Text "robleme." "robleme"
Replace-Text "Effizienter Ressourcenverbrauch." "Effizienter Ressourcenverbrauch"
Replace-Text "Ist für jede Person leicht bedienbar. " "Ist für jede Person leicht bedienbar"

# Globale Testfälle section
Replace-Text "In öffentlichen Netzwerken testen, ob online Modus funktioniert." "In öffentlichen Netzwerken testen, ob online Modus funktioniert"

# Entwicklungsumgebung section
Replace-Text "Android Studio." "Android Studio"
